# daily auto push: 2026-02-28 06:55 UTC
# Insert one new data row right after the existing "2026/02/28" rows (row 891),
# pushing every subsequent row down by one. The new row records another
# observation (hour 13) for 2026/02/28 (a Saturday -> "土").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 892:933 down to 893:934 before writing the new row, so we don't
# clobber data we still need to copy.
$ws.Rows.Item(892).EntireRow.Insert()

$ws.Cells.Item(892, 1).Value = "2026/02/28"
$ws.Cells.Item(892, 2).Value = "土"
$ws.Cells.Item(892, 3).Value = 13
$ws.Cells.Item(892, 4).Value = 201
